$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parametrosInicio")
$ws.Activate()

# Fecha (B5) and Periodo (B6) values filled in by the bot GUI
$ws.Range("B5").Value = "31.01.2023"
$ws.Range("B6").Value = 10

# Leave the selection where the user last clicked before saving
$ws.Range("E7").Select()
